$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column C, shifting ExpPoints from C to G
$ws.Range("C1:F1").EntireColumn.Insert()

# Copy the header style (bold, border, centered) from A1 onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)

# New header cell text
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# Fill C2:F19 with empty strings (inline string cells) for each data row
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
}
